$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.709.15"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "1.635.27"
$ws.Range("E3").Value = "  -0.61%  "
$ws.Range("E4").Value = "  +0.38%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.93"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.503"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.81%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  -0.81%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0621"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.66%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.03"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.59%  "
$ws.Range("E11").Value = "  -0.02%  "
$ws.Range("D12").Value = "1.863.74"
$ws.Range("E12").Value = "  -0.75%  "
$ws.Range("D13").Value = "1.624.08"
$ws.Range("E13").Value = "  -2.32%  "
$ws.Range("E14").Value = "  -0.99%  "
$ws.Range("E15").Value = "  -1.49%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.27"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.28%  "
$ws.Range("D17").Value = "26.677.22"
$ws.Range("E17").Value = "  -0.14%  "
$ws.Range("D18").Value = "0.0₃0725"
$ws.Range("E18").Value = "  -2.34%  "
$ws.Range("E19").Value = "  +0.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "210.41"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -3.56%  "
$ws.Range("E21").Value = "  -0.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.34"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +3.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.17"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.22"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.89%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.32"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.32%  "
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("E27").Value = "  -2.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.04"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.54"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.94%  "
$ws.Range("E30").Value = "  -2.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.18"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.66%  "
$ws.Range("E32").Value = "  -0.47%  "
$ws.Range("E33").Value = "  -1.29%  "
$ws.Range("D34").Value = "1.274.67"
$ws.Range("E34").Value = "  -0.16%  "
$ws.Range("E35").Value = "  -1.26%  "
$ws.Range("E36").Value = "  +0.51%  "
$ws.Range("E37").Value = "  -2.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.530"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.806"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.95%  "
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.800"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.53%  "
$ws.Range("E42").Value = "  -2.22%  "
$ws.Range("D43").Value = "1.774.64"
$ws.Range("E43").Value = "  -0.69%  "
$ws.Range("E44").Value = "  -3.58%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "60.43"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.19"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.66%  "
$ws.Range("E47").Value = "  -2.20%  "
$ws.Range("E48").Value = "  +0.90%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.52"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -3.24%  "
$ws.Range("E50").Value = "  -1.08%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.01"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.16%  "
